# Generate Report for Handoff
# Updates the "Latest Handoff Datetime" (column D, row 4) for the
# 4a773841-d93f-454e-8e9e-9230eb64a205.md file on both the zh-cn and
# de-de localization-status sheets, reflecting a newer handoff report.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D4").Value = "2016-03-03 06:35:25"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D4").Value = "2016-03-03 06:35:36"
